$wb = $excel.ActiveWorkbook

# --- Sheet "Submit orders": add rows 86 and 87 ---
$ws1 = $wb.Worksheets.Item("Submit orders")

$ws1.Range("A86:E87").Style = "Normal"

$ws1.Range("A86").Value = "10.17.2022 12:54 (Kyiv+Israel) 09:54 (UTC) 18:54 (Japan) 15:24 (India)"
$ws1.Range("B86").Value = 1.487
$ws1.Range("C86").Value = -0.7600000000000001
$ws1.Range("D86").Value = "***"
$ws1.Range("E86").Value = "***"

$ws1.Range("A87").Value = "10.17.2022 12:59 (Kyiv+Israel) 09:59 (UTC) 18:59 (Japan) 15:29 (India)"
$ws1.Range("B87").Value = "***"
$ws1.Range("C87").Value = "***"
$ws1.Range("D87").Value = 2.849
$ws1.Range("E87").Value = -1.711

# --- Sheet "Submit a phone survey": add row 76 ---
$ws3 = $wb.Worksheets.Item("Submit a phone survey")

$ws3.Range("A76:E76").Style = "Normal"

$ws3.Range("A76").Value = "10.20.2022 23:26 (Kyiv+Israel) 20:26 (UTC) 05:26 (Japan) 01:56 (India)"
$ws3.Range("B76").Value = 1.707
$ws3.Range("C76").Value = -0.603
$ws3.Range("D76").Value = "***"
$ws3.Range("E76").Value = "***"

# --- Sheet "Checkertificate": add row 91 ---
$ws4 = $wb.Worksheets.Item("Checkertificate")

$ws4.Range("A91:E91").Style = "Normal"

$ws4.Range("A91").Value = "10.20.2022 23:35 (Kyiv+Israel) 20:35 (UTC) 05:35 (Japan) 02:05 (India)"
$ws4.Range("B91").Value = 0.849
$ws4.Range("C91").Value = -0.1839999999999999
$ws4.Range("D91").Value = "***"
$ws4.Range("E91").Value = "***"
